$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9, P9: price text "6.0000" -> "20.0000". The cell's number format
# ("0.00") makes a plain .Value assignment of a numeric-looking string get
# auto-coerced into a real number, which would change the stored cell type.
# Temporarily flip the format to text, write the literal string, then put
# the original numeric display format back so the style stays identical.
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "20.0000"
$ws.Range("P9").NumberFormat = "0.00"

# Row 9, Q9: transaction-count text "3:0" -> "10:0" (already plain text,
# non-numeric-looking, so no coercion risk).
$ws.Range("Q9").Value = "10:0"

# Row 10, P10: running total bumped from 93 -> 107 (genuine number cell).
$ws.Range("P10").Value = 107

# Row 11, A11: refreshed export timestamp.
$ws.Range("A11").Value = "Wednesday, 6 August, 2025 9:28 AM"
